# PlayerPerformance_6038.xlsx edit
# - Insert a new "Player Info" worksheet before the existing "ODI Batting" sheet,
#   carrying basic player identity data (ID/NAME/BATTING_HAND/BOWL_STYLE).
# - On "ODI Batting", rename column D from MATCH_CARD_LINK to MATCH_CODE and
#   replace each full howstat.com scorecard URL with just the trailing match
#   code number that used to be the `MatchCode=` query parameter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ODI Batting" is the sheet already in the workbook.
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 2. Insert the new "Player Info" sheet immediately before "ODI Batting" so it
#    becomes the first (left-most) tab, matching the target tab order.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

# NOTE: `$batting` was captured via a 1-based *positional* lookup, and that
# position now refers to the newly-inserted "Player Info" sheet. Re-resolve
# "ODI Batting" by name so the rest of the script edits the right sheet.
$batting = $wb.Worksheets.Item("ODI Batting")

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row
$playerInfo.Range("A2").Value = "6038"
$playerInfo.Range("B2").Value = "Philip Dean Salt"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# Style the header row like the header row already used on "ODI Batting"
# (bold text, thin box border, centered horizontally, top-aligned vertically).
$playerInfoHeader = $playerInfo.Range("A1:D1")
$playerInfoHeader.Font.Bold = $true
$playerInfoHeader.HorizontalAlignment = -4108
$playerInfoHeader.VerticalAlignment = -4160
$playerInfoHeader.Borders.LineStyle = 1
$playerInfoHeader.Borders.Weight = 2

$playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and shrink each full
#    scorecard URL down to the bare match-code number.
# ---------------------------------------------------------------------------
$batting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("4472", "4473", "4476", "4598", "4599", "4602", "4620", "4622", "4660", "4663", "4666", "4711", "4713", "4717")
for ($i = 0; $i -lt $matchCodes.Count; $i++) {
    $row = $i + 2
    $batting.Cells.Item($row, 4).Value = $matchCodes[$i]
}
